# Fruta / hortaliza, semanal
# Inserts the latest weekly price observation for Ajo (Femacal de La Calera)
# as a new row 239, pushing the previously-existing rows 239-247 down to
# 240-248 (dimension grows from A1:R247 to A1:R248).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 239, shifting rows 239:247
# down to 240:248.
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with this week's data point.
$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44509
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = 100112003
$ws.Cells.Item(239, 7).Value = "Ajo"
$ws.Cells.Item(239, 8).Value = "Chino"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 70
$ws.Cells.Item(239, 11).Value = 16000
$ws.Cells.Item(239, 12).Value = 16500
$ws.Cells.Item(239, 13).Value = 16250
$ws.Cells.Item(239, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(239, 15).Value = "China"
$ws.Cells.Item(239, 16).Value = 1625
$ws.Cells.Item(239, 17).Value = 10
$ws.Cells.Item(239, 18).Value = "Hortaliza"
